$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4056753516197205
$ws.Range("B1").Value = 1.357993483543396
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.001107692718506
$ws.Range("E1").Value = 1.136234283447266
